# Append two new test-case rows (t8, t9) to the table on Sheet1,
# mirroring the existing Testcase/t3..t7 rows in columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "t8"
$ws.Range("B9").Value = "e"
$ws.Range("C9").Value = "y"
$ws.Range("D9").Value = "k"

# Row 10 (column B written before A so the shared-string table order
# matches what Excel produced for this edit: t8, k, o, t9, m)
$ws.Range("B10").Value = "o"
$ws.Range("A10").Value = "t9"
$ws.Range("C10").Value = "o"
$ws.Range("D10").Value = "m"

# Leave the selection where Excel would land after entering the data
$ws.Range("D11").Select() | Out-Null
